$p = $ppt.ActivePresentation

# The deck's current last slide ("Merci pour votre temps") becomes the
# second-to-last slide; a duplicate of it is appended as the new last
# slide, keeping the closing "thank you" slide intact at the very end
# while the (formerly) last slide is repurposed into a transition slide
# introducing the demo section.
$lastIndex = $p.Slides.Count
$lastSlide = $p.Slides.Item($lastIndex)
$lastSlide.Duplicate() | Out-Null

# Repurpose the original last slide into the "moving on to demos" slide.
$s = $p.Slides.Item($lastIndex)
$shp = $s.Shapes.Item(6)
$shp.TextFrame.TextRange.Text = "On passe maintenant pour voire les démos ☺"
$shp.Height = 92.09063720703125
